$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated "Datos actualizados" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 09:22"

# --- Pure numeric stat updates (country unchanged) ---
$ws.Range("B17").Value = 12386
$ws.Range("C17").Value = 89
$ws.Range("E17").Value = 8703

$ws.Range("B39").Value = 2752
$ws.Range("C39").Value = 147
$ws.Range("E39").Value = 2163

$ws.Range("D74").Value = 47
$ws.Range("E74").Value = 617

# --- Country re-labels with updated stats (re-sorted country list) ---
$ws.Range("A156").Value = "Gabon"
$ws.Range("B156").Value = 30
$ws.Range("C156").Value = 6
$ws.Range("D156").Value = 1
$ws.Range("E156").Value = 28

$ws.Range("A157").Value = "Benin"
$ws.Range("B157").Value = 26
$ws.Range("D157").Value = 5
$ws.Range("E157").Value = 20

$ws.Range("A158").Value = "Haiti"
$ws.Range("D158").Value = 0
$ws.Range("E158").Value = 23

$ws.Range("A170").Value = "Fiyi"
$ws.Range("C170").Value = 1
$ws.Range("F170").Value = 0

$ws.Range("A171").Value = "Antigua y Barbuda"
$ws.Range("C171").Value = 0
$ws.Range("F171").Value = 1

$ws.Range("A175").Value = "Sudan"
$ws.Range("C175").Value = 2
$ws.Range("D175").Value = 2
$ws.Range("E175").Value = 10
$ws.Range("H175").Value = 2

$ws.Range("A176").Value = "Liberia"
$ws.Range("B176").Value = 14
$ws.Range("D176").Value = 3
$ws.Range("E176").Value = 8
$ws.Range("H176").Value = 3

$ws.Range("A177").Value = "Curazao"
$ws.Range("B177").Value = 13
$ws.Range("D177").Value = 5
$ws.Range("E177").Value = 7
$ws.Range("F177").Value = 0
$ws.Range("H177").Value = 1

$ws.Range("A179").Value = "Granada"
$ws.Range("D179").Value = 0
$ws.Range("E179").Value = 12
$ws.Range("F179").Value = 2
$ws.Range("H179").Value = 0

$ws.Range("A184").Value = "Mozambique"
$ws.Range("D184").Value = 1
$ws.Range("H184").Value = 0

$ws.Range("A185").Value = "Zimbabue"
$ws.Range("D185").Value = 0
$ws.Range("H185").Value = 1

$ws.Range("A193").Value = "Belice"
$ws.Range("D193").Value = 0
$ws.Range("F193").Value = 1
$ws.Range("H193").Value = 1

$ws.Range("A195").Value = "San Vicente y las Granadinas"
$ws.Range("D195").Value = 1
$ws.Range("F195").Value = 0
$ws.Range("H195").Value = 0

$ws.Range("A198").Value = "San Bartolome"
$ws.Range("D198").Value = 1
$ws.Range("H198").Value = 0

$ws.Range("A199").Value = "Nicaragua"
$ws.Range("D199").Value = 0
$ws.Range("H199").Value = 1

# --- Country re-labels only (stats identical, no numeric change needed) ---
$ws.Range("A207").Value = "Burundi"
$ws.Range("A208").Value = "Anguila"
$ws.Range("A209").Value = "Islas Virgenes Britanicas"
$ws.Range("A213").Value = "Timor Oriental"
$ws.Range("A214").Value = "Sudan del Sur"
$ws.Range("A215").Value = "San Pedro y Miquelon"
